$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: add new trailing cells X9, Y9 ---
$ws.Range("X9").Value = -0.21000099999999833
$ws.Range("Y9").Value = "Down"

# --- Row 10: brand new row of data ---
$ws.Range("A10").Value = 42653.879479166666
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Neutral"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = "Random"
$ws.Range("Q10").Value = 47.96375473473072
$ws.Range("R10").Value = 0.49
$ws.Range("S10").Value = 0.0521
$ws.Range("T10").Value = -0.0214
$ws.Range("U10").Value = 2.25
$ws.Range("V10").Value = "N/A"
$ws.Range("W10").Value = 0

# Copy number formats from row 9 (avoids creating new style entries)
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("S9").Copy()
$ws.Range("S10").PasteSpecial(-4122)
$ws.Range("T9").Copy()
$ws.Range("T10").PasteSpecial(-4122)

$excel.CutCopyMode = 0
